$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 header labels
$ws.Range("A2").Value = "in: floor type"
$ws.Range("B2").Value = "in: length"
$ws.Range("C2").Value = "in: width"
$ws.Range("D2").Value = "out: room cost"

# Lookup table (G:H) values
$ws.Range("G2").Value = "hardwood"
$ws.Range("H2").Value = 1.39
$ws.Range("G3").Value = "carpet"
$ws.Range("H3").Value = 3.99
$ws.Range("G4").Value = "tile"
$ws.Range("H4").Value = 4.99

# Data rows 3-7
$ws.Range("A3").Value = "hardwood"
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = 15
$ws.Range("D3").Formula = "=(B3*C3)*H2"

$ws.Range("A4").Value = "carpet"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 15
$ws.Range("D4").Formula = "=(B4*C4)*H3"

$ws.Range("A5").Value = "tile"
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 16
$ws.Range("D5").Formula = "=(B5*C5)*H4"

$ws.Range("A6").Value = "tile"
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 11
$ws.Range("D6").Formula = "=(B6*C6)*H4"

$ws.Range("A7").Value = "hardwood"
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = 10
$ws.Range("D7").Formula = "=(B7*C7)*H2"

# Total row
$ws.Range("D8").Formula = "=SUM(D3:D7)"

# Numeric inputs B3:C7 centered horizontally (establishes the "horizontal-only" style)
$ws.Range("B3:C7").HorizontalAlignment = -4108

# Header A2 centered both ways
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108

# Floor type values A3:A7 centered both ways (reuses B3:C7's horizontal style, then adds vertical)
$ws.Range("A3:A7").HorizontalAlignment = -4108
$ws.Range("A3:A7").VerticalAlignment = -4108

# Header D2 centered both ways
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").VerticalAlignment = -4108

# Cost values D3:D8: centered both ways (reuses A3:A7's style), then currency format
$ws.Range("D3:D8").HorizontalAlignment = -4108
$ws.Range("D3:D8").VerticalAlignment = -4108
$ws.Range("D3:D8").NumberFormat = "_([`$`$-409]* #,##0.00_);_([`$`$-409]* \(#,##0.00\);_([`$`$-409]* ""-""??_);_(@_)"

$ws.Range("B3").Select() | Out-Null
